# The commit swaps the contents of ppt/theme/theme1.xml and
# ppt/theme/theme2.xml: theme2.xml (the theme actually used by the
# deck's slide master / the visible design) changes from the colourful
# "Integral" / "Red Violet" palette to the plain "Office Theme" / "Office"
# palette (and theme1.xml - used only by the Notes Master and otherwise
# inert - gets the old Integral content; that part isn't reachable
# through the PowerPoint object model, so we reproduce the
# user-visible, design-driving half of the swap: the live colour
# scheme of the presentation's design/master).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$cs = $s.ColorScheme

# New "Office Theme" colour scheme (was previously in theme1.xml,
# becomes the content of theme2.xml -- the theme bound to the slide
# master that actually renders the deck).
$cs.Colors(1).RGB  = 0          # dk1      000000
$cs.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$cs.Colors(3).RGB  = 6968388    # dk2      44546A
$cs.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$cs.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$cs.Colors(6).RGB  = 3243501    # accent2  ED7D31
$cs.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$cs.Colors(8).RGB  = 49407      # accent4  FFC000
$cs.Colors(9).RGB  = 12874308   # accent5  4472C4
$cs.Colors(10).RGB = 4697456    # accent6  70AD47
$cs.Colors(11).RGB = 12673797   # hlink    0563C1
$cs.Colors(12).RGB = 7491477    # folHlink 954F72
